$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("D3").Value = "Neodymium"
    $ws.Range("E3").Value = "Dysprosium"
    $ws.Range("F3").Value = "Copper ores and concentrates"
    $ws.Range("G3").Value = "Raw silicon"
}
